$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2420.9167
$ws.Range("J70").Value = 1814.4286
$ws.Range("L70").Value = 5443.2858
$ws.Range("N70").Value = -5983.2858
$ws.Range("H73").Value = 2420.9167
$ws.Range("J73").Value = 1814.4286
$ws.Range("L73").Value = 5443.2858
$ws.Range("N73").Value = -7315.2858
$ws.Range("H96").Value = 779.9
$ws.Range("J96").Value = 349.5
$ws.Range("L96").Value = 1048.5
$ws.Range("N96").Value = -3794.5
$ws.Range("H103").Value = 526.6
$ws.Range("I103").Value = 211.33333
$ws.Range("J103").Value = 999.5
$ws.Range("K103").Value = 633.99999
$ws.Range("L103").Value = 2998.5
$ws.Range("M103").Value = -47.99999000000003
$ws.Range("N103").Value = -4170.5
$ws.Range("H113").Value = 6983.8335
$ws.Range("I113").Value = 6976.25
$ws.Range("K113").Value = 6976.25
$ws.Range("M113").Value = -3722.25
$ws.Range("H138").Value = 4402.635
$ws.Range("I138").Value = 4146.871
$ws.Range("K138").Value = 12440.613
$ws.Range("M138").Value = -7300.613000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3059.375
$ws.Range("I61").Value = 1650
$ws.Range("K61").Value = 1650
$ws.Range("M61").Value = -1438
$ws.Range("H74").Value = 14290330
$ws.Range("I74").Value = 22728960
$ws.Range("K74").Value = 22728960
$ws.Range("M74").Value = -22728086
$ws.Range("H77").Value = 14290330
$ws.Range("I77").Value = 22728960
$ws.Range("K77").Value = 113644800
$ws.Range("M77").Value = -113640432
$ws.Range("H88").Value = 2913.9167
$ws.Range("I88").Value = 2790
$ws.Range("J88").Value = 2955.2222
$ws.Range("K88").Value = 2790
$ws.Range("L88").Value = 2955.2222
$ws.Range("M88").Value = -2384
$ws.Range("N88").Value = -3767.2222
$ws.Range("H91").Value = 2913.9167
$ws.Range("I91").Value = 2790
$ws.Range("J91").Value = 2955.2222
$ws.Range("K91").Value = 2790
$ws.Range("L91").Value = 2955.2222
$ws.Range("M91").Value = -1386
$ws.Range("N91").Value = -5763.2222
$ws.Range("H136").Value = 3059.375
$ws.Range("I136").Value = 1650
$ws.Range("K136").Value = 4950
$ws.Range("M136").Value = -2400
$ws.Range("H141").Value = 16214
$ws.Range("I141").Value = 12000
$ws.Range("J141").Value = 20428
$ws.Range("K141").Value = 12000
$ws.Range("L141").Value = 20428
$ws.Range("M141").Value = -6820
$ws.Range("N141").Value = -30788

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1071.6296
$ws.Range("I20").Value = 913.75
$ws.Range("J20").Value = 1301.2727
$ws.Range("K20").Value = 913.75
$ws.Range("L20").Value = 1301.2727
$ws.Range("M20").Value = -666.75
$ws.Range("N20").Value = -1795.2727
$ws.Range("H33").Value = 8500
$ws.Range("I33").Value = 7000
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 7000
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -6664
$ws.Range("N33").Value = -10672
$ws.Range("H64").Value = 1192.6666
$ws.Range("I64").Value = 1386.6666
$ws.Range("K64").Value = 1386.6666
$ws.Range("M64").Value = -1161.6666
$ws.Range("H67").Value = 1192.6666
$ws.Range("I67").Value = 1386.6666
$ws.Range("K67").Value = 1386.6666
$ws.Range("M67").Value = -606.6666
$ws.Range("H86").Value = 2091.625
$ws.Range("I86").Value = 2442.375
$ws.Range("J86").Value = 1740.875
$ws.Range("K86").Value = 2442.375
$ws.Range("L86").Value = 1740.875
$ws.Range("M86").Value = -1319.375
$ws.Range("N86").Value = -3986.875
$ws.Range("H89").Value = 2091.625
$ws.Range("I89").Value = 2442.375
$ws.Range("J89").Value = 1740.875
$ws.Range("K89").Value = 2442.375
$ws.Range("L89").Value = 1740.875
$ws.Range("M89").Value = -6595.875
$ws.Range("N89").Value = -19936.375
$ws.Range("H95").Value = 48887.668
$ws.Range("J95").Value = 48887.668
$ws.Range("L95").Value = 48887.668
$ws.Range("N95").Value = -54379.668
$ws.Range("H105").Value = 3610.158
$ws.Range("I105").Value = 3544.0588
$ws.Range("J105").Value = 4172
$ws.Range("K105").Value = 3544.0588
$ws.Range("L105").Value = 4172
$ws.Range("M105").Value = -1797.0588
$ws.Range("N105").Value = -7666
$ws.Range("H107").Value = 2603.8572
$ws.Range("I107").Value = 2635.5454
$ws.Range("K107").Value = 2635.5454
$ws.Range("M107").Value = -715.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3518.3333
$ws.Range("I58").Value = 2055.3333
$ws.Range("K58").Value = 2055.3333
$ws.Range("M58").Value = -1852.3333
$ws.Range("H99").Value = 2067.0652
$ws.Range("I99").Value = 1992.4147
$ws.Range("K99").Value = 1992.4147
$ws.Range("M99").Value = -494.4147
$ws.Range("H122").Value = 3452.7585
$ws.Range("I122").Value = 3401.8948
$ws.Range("J122").Value = 3549.4
$ws.Range("K122").Value = 10205.6844
$ws.Range("L122").Value = 10648.2
$ws.Range("M122").Value = -7755.6844
$ws.Range("N122").Value = -15548.2
$ws.Range("H126").Value = 2067.0652
$ws.Range("I126").Value = 1992.4147
$ws.Range("K126").Value = 5977.2441
$ws.Range("M126").Value = -3507.2441
$ws.Range("H132").Value = 5624.3335
$ws.Range("I132").Value = 3693.6875
$ws.Range("J132").Value = 11802.4
$ws.Range("K132").Value = 11081.0625
$ws.Range("L132").Value = 35407.2
$ws.Range("M132").Value = -8551.0625
$ws.Range("N132").Value = -40467.2
$ws.Range("H136").Value = 3518.3333
$ws.Range("I136").Value = 2055.3333
$ws.Range("K136").Value = 6165.999899999999
$ws.Range("M136").Value = -3615.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 7787.5
$ws.Range("I14").Value = 7787.5
$ws.Range("K14").Value = 23362.5
$ws.Range("M14").Value = -23189.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 650.4783
$ws.Range("I97").Value = 499.1875
$ws.Range("J97").Value = 996.2857
$ws.Range("K97").Value = 499.1875
$ws.Range("L97").Value = 996.2857
$ws.Range("M97").Value = -3.1875
$ws.Range("N97").Value = -1988.2857
$ws.Range("H102").Value = 2907.8333
$ws.Range("I102").Value = 1188.25
$ws.Range("J102").Value = 6347
$ws.Range("K102").Value = 1188.25
$ws.Range("L102").Value = 6347
$ws.Range("M102").Value = 433.75
$ws.Range("N102").Value = -9591

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 105320.2
$ws.Range("I7").Value = 131025.5
$ws.Range("J7").Value = 2499
$ws.Range("K7").Value = 131025.5
$ws.Range("L7").Value = 2499
$ws.Range("M7").Value = -130913.5
$ws.Range("N7").Value = -2723
$ws.Range("H40").Value = 16531.182
$ws.Range("I40").Value = 6478.4287
$ws.Range("K40").Value = 6478.4287
$ws.Range("M40").Value = -6342.4287
$ws.Range("H46").Value = 2019.8334
$ws.Range("I46").Value = 760.2381
$ws.Range("J46").Value = 4958.8887
$ws.Range("K46").Value = 760.2381
$ws.Range("L46").Value = 4958.8887
$ws.Range("M46").Value = -572.2381
$ws.Range("N46").Value = -5334.8887
$ws.Range("H68").Value = 7037.3335
$ws.Range("I68").Value = 7605.3335
$ws.Range("J68").Value = 5333.3335
$ws.Range("K68").Value = 7605.3335
$ws.Range("L68").Value = 5333.3335
$ws.Range("M68").Value = -6856.3335
$ws.Range("N68").Value = -6831.3335
$ws.Range("H71").Value = 7037.3335
$ws.Range("I71").Value = 7605.3335
$ws.Range("J71").Value = 5333.3335
$ws.Range("K71").Value = 38026.6675
$ws.Range("L71").Value = 26666.6675
$ws.Range("M71").Value = -34282.6675
$ws.Range("N71").Value = -34154.6675
$ws.Range("H100").Value = 3720
$ws.Range("I100").Value = 3975
$ws.Range("K100").Value = 3975
$ws.Range("M100").Value = -3434
$ws.Range("H122").Value = 10615.333
$ws.Range("I122").Value = 3999.5
$ws.Range("K122").Value = 11998.5
$ws.Range("M122").Value = -9548.5
$ws.Range("H126").Value = 105320.2
$ws.Range("I126").Value = 131025.5
$ws.Range("J126").Value = 2499
$ws.Range("K126").Value = 393076.5
$ws.Range("L126").Value = 7497
$ws.Range("M126").Value = -390606.5
$ws.Range("N126").Value = -12437

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3649.1
$ws.Range("I96").Value = 3833.3333
$ws.Range("J96").Value = 3570.1428
$ws.Range("K96").Value = 3833.3333
$ws.Range("L96").Value = 3570.1428
$ws.Range("M96").Value = -2460.3333
$ws.Range("N96").Value = -6316.1428
$ws.Range("H122").Value = 2566.5264
$ws.Range("I122").Value = 2542.5
$ws.Range("K122").Value = 7627.5
$ws.Range("M122").Value = -5177.5
$ws.Range("H126").Value = 2738.7334
$ws.Range("I126").Value = 2652.5386
$ws.Range("K126").Value = 7957.6158
$ws.Range("M126").Value = -5487.6158
